$p = $ppt.ActivePresentation
$p.AddTitleMaster()
